$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1: replace the old ad-hoc test data with the new "fill all
# mandatory fields" automation step table.
# ---------------------------------------------------------------
$ws1 = $wb.Sheets.Item(1)

# Wipe out the old contents (A1:E10 covers the previous used range) so
# stale shared strings get dropped on save.
$ws1.Range("A1:E10").ClearContents()

# B2 is the only cell in the old layout that keeps a leftover "highlight"
# style (s="1") that must not survive onto the new header/Action column -
# reset it back to the default look explicitly.
$ws1.Range("B2").Style = "Normal"

# Header row
$ws1.Range("A1").Value = "Feature Step Name"
$ws1.Range("B1").Value = "Action"
$ws1.Range("C1").Value = "Xpath"
$ws1.Range("D1").Value = "Value"

# Row 2 - first name
$ws1.Range("A2").Value = "fill_all_mandatory_fields"
$ws1.Range("B2").Value = "Sendkeys"
$ws1.Range("C2").Value = '//*[@id=\"firstname\"]'
$ws1.Range("D2").Value = "Automation 1"

# Row 3 - last name
$ws1.Range("A3").Value = "fill_all_mandatory_fields"
$ws1.Range("B3").Value = "Sendkeys"
$ws1.Range("C3").Value = '//*[@id="lastname"]'
$ws1.Range("D3").Value = "test"

# Row 4 - email address
$ws1.Range("A4").Value = "fill_all_mandatory_fields"
$ws1.Range("B4").Value = "Sendkeys"
$ws1.Range("C4").Value = '//*[@id="email_address"]'
$ws1.Range("D4").Value = "emailmaster@mailinator.com"

# Row 5 - password
$ws1.Range("A5").Value = "fill_all_mandatory_fields"
$ws1.Range("B5").Value = "Sendkeys"
$ws1.Range("C5").Value = '//*[@id="password"]'
$ws1.Range("D5").Value = "Test@123"

# Row 6 - password confirmation
$ws1.Range("A6").Value = "fill_all_mandatory_fields"
$ws1.Range("B6").Value = "Sendkeys"
$ws1.Range("C6").Value = '//*[@id="password-confirmation"]'
$ws1.Range("D6").Value = "Test@123"

# Row 7 - gender
$ws1.Range("A7").Value = "fill_all_mandatory_fields"
$ws1.Range("B7").Value = "Select"
$ws1.Range("C7").Value = '//*[@id=\"gender\"]'
$ws1.Range("D7").Value = "Male"

# Re-apply the highlighted (Consolas, teal) font style used on the
# "Feature Step Name"/"Xpath"/"Value" data columns (A, C, D) for rows 2-7.
# Style a single "template" cell directly (this is the only font mutation
# that cleanly reuses the workbook's existing named font/style), then
# fan the formatting out to every other cell with a Copy/PasteSpecial
# (formats only) so the shared cell values are left untouched.
$f = $ws1.Range("A2").Font
$f.Color = 10733079
$f.Name = "Consolas"
$f.Size = 10

$ws1.Range("A2").Copy()
$ws1.Range("A3:A7").PasteSpecial(-4122)
$ws1.Range("C2:C7").PasteSpecial(-4122)
$ws1.Range("D2:D7").PasteSpecial(-4122)

# Column widths (A, B, C, D)
$ws1.Columns.Item(1).ColumnWidth = 29.3833
$ws1.Columns.Item(2).ColumnWidth = 66.3833
$ws1.Columns.Item(3).ColumnWidth = 46.7167
$ws1.Columns.Item(4).ColumnWidth = 34.7167

# Update selection to match the saved view state
$ws1.Range("B6").Select()

# ---------------------------------------------------------------
# Sheet2: content/text is unchanged, only the underlying shared
# string table shrinks as a side effect of removing the old Sheet1
# strings above - nothing further to do here.
# ---------------------------------------------------------------
